# NIT-9003321466.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The worker table used to hold 3 rows: MARILUZ TORO ORTEGA (row 16),
# PEDRO PATERNINA CARBALLO (row 17) and a totals/closing row for the NIT
# itself (row 18). The update removes the obsolete NIT closing row and
# keeps only the two real workers, now with PEDRO listed first (row 16)
# and MARILUZ second (row 17, which becomes the new last row of the
# table and therefore inherits the solid bottom border that used to
# close the table at row 18). The "Cant. Trabajadores" / "Cant.
# Periodos" counters drop from 3 to 2, and the "VALOR MORA" total is
# refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old NIT totals row (9003321466 / 1606 / 91120 / 0) — the
# two worker rows below it shift up automatically (old 23/24 -> 22/23).
$ws.Rows("18:18").Delete()

# Row 16 now becomes PEDRO PATERNINA CARBALLO's record.
$ws.Range("C16").Value = "73200842"
$ws.Range("D16").Value = "PEDRO PATERNINA CARBALLO"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 27578

# Row 17 now becomes MARILUZ TORO ORTEGA's record.
$ws.Range("C17").Value = "45690981"
$ws.Range("D17").Value = "MARILUZ TORO ORTEGA"
$ws.Range("E17").Value = "1704"
$ws.Range("F17").Value = 15628

# Row 17 is now the last row of the table, so it gets the solid black
# bottom border that used to close the table at the deleted row 18.
$closingRow = $ws.Range("B17:J17")
$closingRow.Borders.Item(9).LineStyle = 1
$closingRow.Borders.Item(9).Color = 0

# Refresh the summary figures for the now-2-worker statement.
$ws.Range("E11").Value = 43206
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
